# Update data: 8 April 2022
# Adds the new March-2022 (Excel serial date 44621) observations to both
# the "Canada" sheet (sheet1) and the "Province" sheet (sheet2).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Canada" (sheet1): append row 28 for Canada, March 2022.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Canada")

# Copy the formatting of the previous data row down one row so the new
# row picks up the same date-number-format style used throughout the
# column (matches style index reused by the source workbook).
$ws1.Range("A27:E27").Copy()
$ws1.Range("A28").PasteSpecial(-4122)

$ws1.Range("A28").Value = 44621
$ws1.Range("B28").Value = "Canada"
$ws1.Range("D28").Value = 1100.2
$ws1.Range("E28").Value = 1169.2
$ws1.Range("C28").Formula = "=(D28-E28)/E28*100"

$ws1.Range("A28").Select()

# ---------------------------------------------------------------------
# Sheet "Province" (sheet2): append rows 262-271, one per province, for
# March 2022 (same date serial 44621).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Province")

# Copy the formatting of the prior date-block (rows 252-261) down into
# the new block (rows 262-271) so formatting/style indices match.
$ws2.Range("A252:E261").Copy()
$ws2.Range("A262").PasteSpecial(-4122)

$ws2.Range("A262").Value = 44621
$ws2.Range("B262").Value = "Newfoundland & Labrador"
$ws2.Range("C262").Formula = "=(D262-E262)/E262*100"
$ws2.Range("D262").Value = 33.3
$ws2.Range("E262").Value = 31.6

$ws2.Range("A263").Value = 44621
$ws2.Range("B263").Value = "Prince Edward Island"
$ws2.Range("C263").Formula = "=(D263-E263)/E263*100"
$ws2.Range("D263").Value = 7.5
$ws2.Range("E263").Value = 7.6

$ws2.Range("A264").Value = 44621
$ws2.Range("B264").Value = "Nova Scotia"
$ws2.Range("C264").Formula = "=(D264-E264)/E264*100"
$ws2.Range("D264").Value = 32.8
$ws2.Range("E264").Value = 32

$ws2.Range("A265").Value = 44621
$ws2.Range("B265").Value = "New Brunswick"
$ws2.Range("C265").Formula = "=(D265-E265)/E265*100"
$ws2.Range("D265").Value = 30.4
$ws2.Range("E265").Value = 31.6

$ws2.Range("A266").Value = 44621
$ws2.Range("B266").Value = "Quebec"
$ws2.Range("C266").Formula = "=(D266-E266)/E266*100"
$ws2.Range("D266").Value = 190
$ws2.Range("E266").Value = 238.5

$ws2.Range("A267").Value = 44621
$ws2.Range("B267").Value = "Ontario"
$ws2.Range("C267").Formula = "=(D267-E267)/E267*100"
$ws2.Range("D267").Value = 429.8
$ws2.Range("E267").Value = 462.4

$ws2.Range("A268").Value = 44621
$ws2.Range("B268").Value = "Manitoba"
$ws2.Range("C268").Formula = "=(D268-E268)/E268*100"
$ws2.Range("D268").Value = 37
$ws2.Range("E268").Value = 35.8

$ws2.Range("A269").Value = 44621
$ws2.Range("B269").Value = "Saskatchewan"
$ws2.Range("C269").Formula = "=(D269-E269)/E269*100"
$ws2.Range("D269").Value = 30.5
$ws2.Range("E269").Value = 30.6

$ws2.Range("A270").Value = 44621
$ws2.Range("B270").Value = "Alberta"
$ws2.Range("C270").Formula = "=(D270-E270)/E270*100"
$ws2.Range("D270").Value = 162.5
$ws2.Range("E270").Value = 166.9

$ws2.Range("A271").Value = 44621
$ws2.Range("B271").Value = "British Columbia"
$ws2.Range("C271").Formula = "=(D271-E271)/E271*100"
$ws2.Range("D271").Value = 146.4
$ws2.Range("E271").Value = 132.3

# Selecting D272 last (matches the post-edit selection in the source
# workbook) also restores the "Province" tab as the active sheet tab.
$ws2.Range("D272").Select()
